# rill-analysis: compatible enhancement for text/plain content-type but
# json content response.
#
# The "_input" sheet gains a new row ("数据粒度" / "按月查看") inserted right
# above the former row 4 ("时间" / date). Everything that referenced the old
# row 4 / row 5 positions (the trend chart series, the defined name used to
# build the chart's data source range, and the report title formula) is
# updated to point one row further down.

$wb = $excel.ActiveWorkbook

$inputWs = $wb.Worksheets.Item("_input")
$trendWs = $wb.Worksheets.Item("trend")

# --- 1. Insert the new "数据粒度 / 按月查看" row above the old row 4 -------------
$inputWs.Rows("4:4").Insert()
$inputWs.Rows("4:4").RowHeight = 18

$inputWs.Range("A4").Value2 = "数据粒度"
$inputWs.Range("B4").Value2 = "按月查看"

# --- 2. Update the workbook-level defined name used as the chart source ----
$dataSourceName = $wb.Names.Item("trendChartDataSource")
$dataSourceName.RefersTo = "=OFFSET(_input!`$A5,0,0,COUNTA(_input!`$A:`$A)-4, COUNTA(_input!`$5:`$5))"

# --- 3. Update the trend chart series to point at the shifted rows ---------
$chartObj = $trendWs.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(_input!`$A`$6,_input!`$B`$5:`$B`$5,_input!`$B`$6:`$B`$6,1)"

# --- 4. Update the report title formula on the trend sheet -----------------
$trendWs.Range("B2").Formula = "=_input!`$B2&_input!`$B3&_input!`$B4&`"报表`""
